$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove existing hyperlinks first; row insertion below does not
#    automatically re-bind hyperlink objects to their shifted cells in
#    this runtime, so we recreate them afterwards at the correct spots.
$ws.Cells.Hyperlinks.Delete()

# 2. Insert the new row for the "Focal nodular hyperplasia - Isoechoic"
#    clip, pushing every row from 11 downward by one.
$ws.Rows.Item(11).Insert()

# 3. Populate the newly inserted row 11.
$ws.Range("A11").Value = "Liver"
$ws.Range("B11").Value = "Focal nodular hyperplasia - Isoechoic "
$ws.Range("C11").Value = "Clip 2 B-mode + Color + microV"
$ws.Range("D11").Value = "https://youtu.be/PGVchRMB22g"

# 4. Refresh the stored sort-state range/condition to account for the
#    extra row (was A2:C26 / A2:A26, becomes A2:C27 / A2:A27). We first
#    apply across A:D so the already-sorted rows stay aligned, then
#    restate (and re-apply, a no-op on ordering) the narrower A:C range
#    that matches the original sort definition.
$ws.Sort.SetRange($ws.Range("A2:D27"))
$ws.Sort.SortFields.Add($ws.Range("A2:A27")) | Out-Null
$ws.Sort.Apply()
$ws.Sort.SetRange($ws.Range("A2:C27"))
$ws.Sort.Apply()

# 5. Recreate the hyperlinks at their final (post-insert) cell locations.
#    Hyperlinks.Add() stamps the anchor cell with a freshly-cloned
#    "hyperlink" cell format; restate the original named style right
#    after each call so the cell keeps referencing the workbook's
#    existing "Collegamento ipertestuale" format instead of the clone.
$linkCells = @("D3","D4","D6","D9","D12","D14","D18","D19","D23","D24","D25","D26","D27","D28","D29")
$linkUrls = @(
  "https://youtu.be/zxTC0YBY2RY",
  "https://youtu.be/K2Wbg7BgXy4",
  "https://youtu.be/2kRZcpi70Aw",
  "https://youtu.be/91M82AIMyu0",
  "https://youtu.be/15o_Km86IzM",
  "https://youtu.be/RhSUFLTmTl4",
  "https://youtu.be/DjI1kEnzfSQ",
  "https://youtu.be/U3ydTsRwxok",
  "https://youtu.be/xBfd04F4Ni8",
  "https://youtu.be/JvwODCASLYQ",
  "https://youtu.be/pc-vbxSRTbs",
  "https://youtu.be/Axbee4vjNtU",
  "https://youtu.be/qushjTAy6XQ",
  "https://youtu.be/_FckFwJwynI",
  "https://youtu.be/z_oaRVxRz5s"
)
for ($i = 0; $i -lt $linkCells.Length; $i++) {
  $ws.Hyperlinks.Add($ws.Range($linkCells[$i]), $linkUrls[$i]) | Out-Null
  $ws.Range($linkCells[$i]).Style = "Collegamento ipertestuale"
}

# 6. Update the sheet view: scroll so column B is the left-most visible
#    column, and move the active selection to D11 (the new row).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D11").Select()
